$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 'Gesellschaften' wurden in 'Firmen'/'Unternehmen' umbenannt
$ws.Range("A2").Value = "Unternehmen"
$ws.Range("A3").Value = "Unternehmenskuerzel"

# Auswahl auf A5 setzen
$ws.Range("A5").Select()
